# Change cell B11 on the "Rules" sheet from the text "R40" to the text "1".
#
# Note: B11's number format is General, and simply assigning a numeric-
# looking string (Value = "1") would make Excel auto-convert it to a
# *number*, while the target content must remain *text* ("1" as a shared
# string) so the cell keeps behaving as the rule-id label column it is
# (R10/R20/R30/"1"), not as a numeric value, and its existing cell style
# (borders/fill) must be left untouched.
#
# The trick Excel itself uses for this is the leading apostrophe ('1),
# which marks the entered text as a literal string ("quote prefix").
# We stage that text in a scratch cell, then copy just the *value* onto
# B11 via PasteSpecial(xlPasteValues) so B11's original formatting/style
# is preserved exactly, and finally clear the scratch cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$scratch.Value = "'1"
$scratch.Copy()

$target = $ws.Range("B11")
$target.PasteSpecial(-4163)   # -4163 = xlPasteValues: paste value+type only, keep destination formatting

$scratch.Clear()
